# Quarterly report update: drop the oldest quarter column (E, "فصل اول منتهی
# به 1399/06") and shift every quarter column one step to the left (F->E,
# G->F, ... N->M), then populate the freed-up rightmost column (N) with the
# newly reported quarter ("فصل سوم منتهی به 1401/12" for the header rows,
# and the new quarter's figures for every data row).
#
# This mirrors how the source workbook is maintained each quarter: the
# 10-quarter rolling window slides forward by one and a new quarter of
# data/prices is appended at the right edge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarter columns E..N (5..14), in order.
$cols = @(5, 6, 7, 8, 9, 10, 11, 12, 13, 14)

# Every row in the sheet that carries quarterly data across E:N.
$rows = @(8, 10, 11, 12, 13, 17, 19, 20, 21, 22, 23, 27, 29, 30, 31, 32, 33, 37, 39, 40, 41, 42, 46, 48, 49, 50, 51, 55, 57, 58, 59, 60)

# New value to drop into column N (last column) for each of those rows,
# once the old contents have been shifted left by one quarter.
$newN = @{}
$newN[8]  = "فصل سوم منتهی به 1401/12"
$newN[17] = "فصل سوم منتهی به 1401/12"
$newN[27] = "فصل سوم منتهی به 1401/12"
$newN[37] = "فصل سوم منتهی به 1401/12"
$newN[46] = "فصل سوم منتهی به 1401/12"
$newN[55] = "فصل سوم منتهی به 1401/12"

$newN[10] = "-"
$newN[19] = "-"
$newN[20] = "-"
$newN[29] = "-"
$newN[30] = "-"
$newN[39] = "-"
$newN[40] = "-"
$newN[41] = "-"
$newN[48] = "-"
$newN[57] = "-"

$newN[11] = 0
$newN[12] = 30285
$newN[13] = 30285
$newN[21] = 0
$newN[22] = 30909
$newN[23] = 30909
$newN[31] = 0
$newN[32] = 3755932
$newN[33] = 3755932
$newN[42] = 121515804
$newN[49] = 0
$newN[50] = -2114339
$newN[51] = -2114339
$newN[58] = 0
$newN[59] = 1641593
$newN[60] = 1641593

foreach ($r in $rows) {
    for ($i = 0; $i -lt ($cols.Length - 1); $i++) {
        $srcCol = $cols[$i + 1]
        $dstCol = $cols[$i]
        $ws.Cells.Item($r, $dstCol).Value = $ws.Cells.Item($r, $srcCol).Value()
    }
    $ws.Cells.Item($r, 14).Value = $newN[$r]
}
